# Update column G ("K" = strikeouts) in sheet1 per regenerated save_data.
# The "K" column previously held a different stat ("Strike#"); values below
# are the recalculated strikeout totals (K) for each start, rows 2-66.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(
    @{ Row = 2; Val = 2 },
    @{ Row = 3; Val = 0 },
    @{ Row = 4; Val = 3 },
    @{ Row = 5; Val = 1 },
    @{ Row = 6; Val = 0 },
    @{ Row = 7; Val = 2 },
    @{ Row = 8; Val = 1 },
    @{ Row = 9; Val = 1 },
    @{ Row = 10; Val = 0 },
    @{ Row = 11; Val = 2 },
    @{ Row = 12; Val = 2 },
    @{ Row = 13; Val = 1 },
    @{ Row = 14; Val = 1 },
    @{ Row = 15; Val = 3 },
    @{ Row = 16; Val = 0 },
    @{ Row = 17; Val = 1 },
    @{ Row = 18; Val = 1 },
    @{ Row = 19; Val = 2 },
    @{ Row = 20; Val = 0 },
    @{ Row = 21; Val = 1 },
    @{ Row = 22; Val = 1 },
    @{ Row = 23; Val = 2 },
    @{ Row = 24; Val = 1 },
    @{ Row = 25; Val = 0 },
    @{ Row = 26; Val = 0 },
    @{ Row = 27; Val = 3 },
    @{ Row = 28; Val = 2 },
    @{ Row = 29; Val = 1 },
    @{ Row = 30; Val = 1 },
    @{ Row = 31; Val = 0 },
    @{ Row = 32; Val = 0 },
    @{ Row = 33; Val = 2 },
    @{ Row = 34; Val = 2 },
    @{ Row = 35; Val = 0 },
    @{ Row = 36; Val = 1 },
    @{ Row = 37; Val = 1 },
    @{ Row = 38; Val = 1 },
    @{ Row = 39; Val = 1 },
    @{ Row = 40; Val = 1 },
    @{ Row = 41; Val = 2 },
    @{ Row = 42; Val = 2 },
    @{ Row = 43; Val = 2 },
    @{ Row = 44; Val = 2 },
    @{ Row = 45; Val = 1 },
    @{ Row = 46; Val = 1 },
    @{ Row = 47; Val = 0 },
    @{ Row = 48; Val = 1 },
    @{ Row = 49; Val = 0 },
    @{ Row = 50; Val = 0 },
    @{ Row = 51; Val = 0 },
    @{ Row = 52; Val = 1 },
    @{ Row = 53; Val = 2 },
    @{ Row = 54; Val = 1 },
    @{ Row = 55; Val = 0 },
    @{ Row = 56; Val = 1 },
    @{ Row = 57; Val = 0 },
    @{ Row = 58; Val = 2 },
    @{ Row = 59; Val = 1 },
    @{ Row = 60; Val = 0 },
    @{ Row = 61; Val = 1 },
    @{ Row = 62; Val = 2 },
    @{ Row = 63; Val = 0 },
    @{ Row = 64; Val = 1 },
    @{ Row = 65; Val = 1 },
    @{ Row = 66; Val = 3 }
)

foreach ($entry in $kValues) {
    $ws.Cells.Item($entry.Row, 7).Value = $entry.Val
}

Write-Output ("Updated {0} K values in column G" -f $kValues.Count)
